$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '30.255.75'
$ws.Range("E2").Value = '  +1.45%  '

# Row 3
$ws.Range("D3").Value = '1.881.77'
$ws.Range("E3").Value = '  +0.20%  '

# Row 4
$ws.Range("D4").Value = '''0.9987'
$ws.Range("E4").Value = '  -0.52%  '

# Row 5
$ws.Range("D5").Value = '''243.81'
$ws.Range("E5").Value = '  -0.08%  '

# Row 6
$ws.Range("D6").Value = '''0.9992'
$ws.Range("E6").Value = '  -0.44%  '

# Row 7
$ws.Range("D7").Value = '''0.4898'
$ws.Range("E7").Value = '  -0.86%  '

# Row 8
$ws.Range("D8").Value = '''0.2909'
$ws.Range("E8").Value = '  +0.21%  '

# Row 9
$ws.Range("D9").Value = '''0.06609'
$ws.Range("E9").Value = '  +0.32%  '

# Row 10
$ws.Range("D10").Value = '1.874.11'
$ws.Range("E10").Value = '  -0.26%  '

# Row 11
$ws.Range("D11").Value = '''16.38'
$ws.Range("E11").Value = '  -2.97%  '

# Row 12
$ws.Range("D12").Value = '''0.07215'
$ws.Range("E12").Value = '  +0.37%  '

# Row 13
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").Value = '''0.6664'
$ws.Range("E13").Value = '  -0.31%  '

# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '''4.972'
$ws.Range("E14").Value = '  +3.31%  '

# Row 15
$ws.Range("D15").Value = '''86.36'
$ws.Range("E15").Value = '  +1.08%  '

# Row 16
$ws.Range("D16").Value = '30.142.27'
$ws.Range("E16").Value = '  +1.00%  '

# Row 17
$ws.Range("D17").Value = '''0.000007802'
$ws.Range("E17").Value = '  -0.02%  '

# Row 18
$ws.Range("D18").Value = '''0.9994'
$ws.Range("E18").Value = '  -0.22%  '

# Row 19
$ws.Range("D19").Value = '''12.80'
$ws.Range("E19").Value = '  +0.34%  '

# Row 20
$ws.Range("D20").Value = '2.116.35'
$ws.Range("E20").Value = '  -0.74%  '

# Row 21
$ws.Range("D21").Value = '''0.9993'
$ws.Range("E21").Value = '  -0.38%  '

# Row 22
$ws.Range("D22").Value = '''4.763'
$ws.Range("E22").Value = '  +0.52%  '

# Row 23
$ws.Range("D23").Value = '''5.919'
$ws.Range("E23").Value = '  +6.08%  '

# Row 24
$ws.Range("D24").Value = '''9.196'
$ws.Range("E24").Value = '  +0.99%  '

# Row 25
$ws.Range("D25").Value = '''152.10'
$ws.Range("E25").Value = '  +2.81%  '

# Row 26
$ws.Range("D26").Value = '''143.99'
$ws.Range("E26").Value = '  +7.30%  '

# Row 27
$ws.Range("D27").Value = '''17.00'
$ws.Range("E27").Value = '  +1.89%  '

# Row 28
$ws.Range("D28").Value = '''1.892'
$ws.Range("E28").Value = '  -1.56%  '

# Row 29
$ws.Range("D29").Value = '''1.400'
$ws.Range("E29").Value = '  +1.54%  '

# Row 30
$ws.Range("D30").Value = '''4.213'
$ws.Range("E30").Value = '  +1.05%  '

# Row 31
$ws.Range("D31").Value = '''0.08814'
$ws.Range("E31").Value = '  +1.60%  '

# Row 32
$ws.Range("D32").Value = '''3.986'
$ws.Range("E32").Value = '  +1.44%  '

# Row 33
$ws.Range("D33").Value = '''0.05225'
$ws.Range("E33").Value = '  +3.18%  '

# Row 34
$ws.Range("D34").Value = '''0.7224'
$ws.Range("E34").Value = '  +2.91%  '

# Row 35
$ws.Range("D35").Value = '''1.113'
$ws.Range("E35").Value = '  +0.50%  '

# Row 36
$ws.Range("E36").Value = '  -0.64%  '

# Row 37
$ws.Range("D37").Value = '''0.01840'
$ws.Range("E37").Value = '  +12.15%  '

# Row 38
$ws.Range("E38").Value = '  -0.37%  '

# Row 39
$ws.Range("D39").Value = '''2.169'
$ws.Range("E39").Value = '  -1.64%  '

# Row 40
$ws.Range("D40").Value = '''0.9313'
$ws.Range("E40").Value = '  -0.15%  '

# Row 41
$ws.Range("D41").Value = '''0.4267'
$ws.Range("E41").Value = '  +2.41%  '

# Row 42
$ws.Range("D42").Value = '''104.14'
$ws.Range("E42").Value = '  +1.49%  '

# Row 43
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '''5.766'
$ws.Range("E43").Value = '  -4.64%  '

# Row 44
$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").Value = '''0.9997'
$ws.Range("E44").Value = '  -0.01%  '

# Row 45
$ws.Range("D45").Value = '''7.443'
$ws.Range("E45").Value = '  +0.04%  '

# Row 46
$ws.Range("D46").Value = '''0.1286'
$ws.Range("E46").Value = '  +2.37%  '

# Row 47
$ws.Range("D47").Value = '''0.05734'
$ws.Range("E47").Value = '  +0.37%  '

# Row 48
$ws.Range("D48").Value = '''32.92'
$ws.Range("E48").Value = '  +1.28%  '

# Row 49
$ws.Range("D49").Value = '''8.322'
$ws.Range("E49").Value = '  +1.71%  '

# Row 50
$ws.Range("D50").Value = '''0.3784'
$ws.Range("E50").Value = '  +2.34%  '

# Row 51
$ws.Range("D51").Value = '''1.354'
$ws.Range("E51").Value = '  +1.38%  '
